$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the step/result values to be numbered per-row so each row has a
# distinct step/result string (fixes import issue where all rows shared the
# same "step"/"result" text). Write all "step N" cells before the "result N"
# cells so the shared string table is built in the same order as the target.
$ws.Range("B1").Value = "step 1"
$ws.Range("B2").Value = "step 2"
$ws.Range("B3").Value = "step 3"

$ws.Range("C1").Value = "result 1"
$ws.Range("C2").Value = "result 2"
$ws.Range("C3").Value = "result 3"

# Move the active selection off the C1:C3 range onto E4, matching the
# post-edit view state.
$ws.Range("E4").Select()
